$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.176.88"
$ws.Range("E2").Value = "  +1.63%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.772.76"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "622.30"
$ws.Range("E5").Value = "  +3.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.14"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.769.83"
$ws.Range("E7").Value = "  -0.45%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("E9").Value = "  +0.46%  "

$ws.Range("E10").Value = "  +1.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.451"
$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.62"
$ws.Range("E12").Value = "  +0.34%  "

$ws.Range("E13").Value = "  -0.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.36"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.408.32"
$ws.Range("E15").Value = "  -0.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.777.52"
$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.169.55"
$ws.Range("E17").Value = "  +1.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.68"
$ws.Range("E18").Value = "  -3.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.08"
$ws.Range("E19").Value = "  +0.40%  "

$ws.Range("E20").Value = "  -1.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "468.17"
$ws.Range("E21").Value = "  +1.62%  "

$ws.Range("E22").Value = "  -0.64%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  +2.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.98"
$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.99"

$ws.Range("E27").Value = "  +2.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.97"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.918.68"
$ws.Range("E30").Value = "  -0.52%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.24"
$ws.Range("E31").Value = "  +0.39%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.65"
$ws.Range("E32").Value = "  +0.76%  "

$ws.Range("E33").Value = "  -0.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.84"
$ws.Range("E34").Value = "  -1.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.722.79"
$ws.Range("E36").Value = "  -0.42%  "

$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.155"
$ws.Range("E38").Value = "  +10.83%  "

$ws.Range("E39").Value = "  +2.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.33"
$ws.Range("E40").Value = "  +0.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.77"
$ws.Range("E41").Value = "  -0.85%  "

$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.966"
$ws.Range("E42").Value = "  -2.18%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "154.20"
$ws.Range("E45").Value = "  +1.37%  "

$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.299"
$ws.Range("E46").Value = "  +0.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "42.87"
$ws.Range("E47").Value = "  -1.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.70"
$ws.Range("E48").Value = "  -1.53%  "

$ws.Range("E49").Value = "  +2.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.39"
$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("E51").Value = "  +0.95%  "
